$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header stays the same (A1="variable", B1="explanation")
$ws.Range("A1").Value = "variable"
$ws.Range("B1").Value = "explanation"

# Rewrite column A for rows 2-63 with the updated/reordered codebook variable list
$ws.Range("A2").Value = "unique_id"
$ws.Range("A3").Value = "exclude_min_max_scores"
$ws.Range("A4").Value = "sub_bfi_a"
$ws.Range("A5").Value = "sub_bfi_c"
$ws.Range("A6").Value = "sub_bfi_e"
$ws.Range("A7").Value = "sub_bfi_o"
$ws.Range("A8").Value = "sub_bfi_n"
$ws.Range("A9").Value = "exclude_iat_completeness"
$ws.Range("A10").Value = "exclude_iat_msseconds"
$ws.Range("A11").Value = "exclude_iat_accuracy"
$ws.Range("A12").Value = "mean_bfi_a"
$ws.Range("A13").Value = "mean_bfi_c"
$ws.Range("A14").Value = "mean_bfi_e"
$ws.Range("A15").Value = "mean_bfi_n"
$ws.Range("A16").Value = "mean_bfi_o"
$ws.Range("A17").Value = "bfi_a1_recode"
$ws.Range("A18").Value = "bfi_a2"
$ws.Range("A19").Value = "bfi_a3_recode"
$ws.Range("A20").Value = "bfi_a4"
$ws.Range("A21").Value = "bfi_a5"
$ws.Range("A22").Value = "bfi_a6_recode"
$ws.Range("A23").Value = "bfi_a7"
$ws.Range("A24").Value = "bfi_a8_recode"
$ws.Range("A25").Value = "bfi_a9"
$ws.Range("A26").Value = "bfi_c1"
$ws.Range("A27").Value = "bfi_c2_recode"
$ws.Range("A28").Value = "bfi_c3"
$ws.Range("A29").Value = "bfi_c4_recode"
$ws.Range("A30").Value = "bfi_c5_recode"
$ws.Range("A31").Value = "bfi_c6"
$ws.Range("A32").Value = "bfi_c7"
$ws.Range("A33").Value = "bfi_c8"
$ws.Range("A34").Value = "bfi_c9_recode"
$ws.Range("A35").Value = "bfi_e1"
$ws.Range("A36").Value = "bfi_e2_recode"
$ws.Range("A37").Value = "bfi_e3"
$ws.Range("A38").Value = "bfi_e4"
$ws.Range("A39").Value = "bfi_e5_recode"
$ws.Range("A40").Value = "bfi_e6"
$ws.Range("A41").Value = "bfi_e7_recode"
$ws.Range("A42").Value = "bfi_e8"
$ws.Range("A43").Value = "bfi_o1"
$ws.Range("A44").Value = "bfi_o2"
$ws.Range("A45").Value = "bfi_o3"
$ws.Range("A46").Value = "bfi_o4"
$ws.Range("A47").Value = "bfi_o5"
$ws.Range("A48").Value = "bfi_o6"
$ws.Range("A49").Value = "bfi_o7_recode"
$ws.Range("A50").Value = "bfi_o8"
$ws.Range("A51").Value = "bfi_o9_recode"
$ws.Range("A52").Value = "bfi_o10"
$ws.Range("A53").Value = "bfi_n1"
$ws.Range("A54").Value = "bfi_n2_recode"
$ws.Range("A55").Value = "bfi_n3"
$ws.Range("A56").Value = "bfi_n4"
$ws.Range("A57").Value = "bfi_n5_recode"
$ws.Range("A58").Value = "bfi_n6"
$ws.Range("A59").Value = "bfi_n7_recode"
$ws.Range("A60").Value = "bfi_n8"
$ws.Range("A61").Value = "age"
$ws.Range("A62").Value = "gender"
$ws.Range("A63").Value = "exclude_participant"

# Drop the now-removed trailing rows 64-78
$ws.Range("A64:B78").Clear()
